# Scheduled runner update: refresh computed market/profit figures across
# several Leve-profit worksheets (ALC, ARM, BSM, CRP, CUL, GSM, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1333.8334
$ws.Range("I70").Value = 1033.3334
$ws.Range("J70").Value = 1634.3334
$ws.Range("K70").Value = 3100.0002
$ws.Range("L70").Value = 4903.0002
$ws.Range("M70").Value = -2830.0002
$ws.Range("N70").Value = -5443.0002
$ws.Range("H73").Value = 1333.8334
$ws.Range("I73").Value = 1033.3334
$ws.Range("J73").Value = 1634.3334
$ws.Range("K73").Value = 3100.0002
$ws.Range("L73").Value = 4903.0002
$ws.Range("M73").Value = -2164.0002
$ws.Range("N73").Value = -6775.0002
$ws.Range("H111").Value = 1261.9286
$ws.Range("I111").Value = 1031.6666
$ws.Range("K111").Value = 3094.9998
$ws.Range("M111").Value = -27.99980000000005
$ws.Range("H125").Value = 1759.5385
$ws.Range("I125").Value = 958.8
$ws.Range("J125").Value = 2260
$ws.Range("K125").Value = 8629.199999999999
$ws.Range("L125").Value = 20340
$ws.Range("M125").Value = -6169.199999999999
$ws.Range("N125").Value = -25260
$ws.Range("H135").Value = 1417.3636
$ws.Range("I135").Value = 1632.6666
$ws.Range("J135").Value = 1159
$ws.Range("K135").Value = 14693.9994
$ws.Range("L135").Value = 10431
$ws.Range("M135").Value = -12158.9994
$ws.Range("N135").Value = -15501
$ws.Range("H137").Value = 2026.7037
$ws.Range("I137").Value = 2184.75
$ws.Range("J137").Value = 1900.2667
$ws.Range("K137").Value = 6554.25
$ws.Range("L137").Value = 5700.800099999999
$ws.Range("M137").Value = -4004.25
$ws.Range("N137").Value = -10800.8001
$ws.Range("H138").Value = 2095.35
$ws.Range("I138").Value = 1083.9333
$ws.Range("J138").Value = 2273.8352
$ws.Range("K138").Value = 3251.7999
$ws.Range("L138").Value = 6821.5056
$ws.Range("M138").Value = 1888.2001
$ws.Range("N138").Value = -17101.5056

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 875590.1
$ws.Range("I32").Value = 1016726.1
$ws.Range("J32").Value = 28774.363
$ws.Range("K32").Value = 1016726.1
$ws.Range("L32").Value = 28774.363
$ws.Range("M32").Value = -1016439.1
$ws.Range("N32").Value = -29348.363
$ws.Range("H88").Value = 3927.182
$ws.Range("I88").Value = 4571.2856
$ws.Range("J88").Value = 2800
$ws.Range("K88").Value = 4571.2856
$ws.Range("L88").Value = 2800
$ws.Range("M88").Value = -4165.2856
$ws.Range("N88").Value = -3612
$ws.Range("H91").Value = 3927.182
$ws.Range("I91").Value = 4571.2856
$ws.Range("J91").Value = 2800
$ws.Range("K91").Value = 4571.2856
$ws.Range("L91").Value = 2800
$ws.Range("M91").Value = -3167.2856
$ws.Range("N91").Value = -5608
$ws.Range("H97").Value = 2391.4285
$ws.Range("I97").Value = 2535
$ws.Range("J97").Value = 2200
$ws.Range("K97").Value = 2535
$ws.Range("L97").Value = 2200
$ws.Range("M97").Value = -2039
$ws.Range("N97").Value = -3192
$ws.Range("H102").Value = 3136.6667
$ws.Range("I102").Value = 3136.6667
$ws.Range("K102").Value = 3136.6667
$ws.Range("M102").Value = -1514.6667
$ws.Range("H122").Value = 93231.82000000001
$ws.Range("I122").Value = 126854.5
$ws.Range("K122").Value = 380563.5
$ws.Range("M122").Value = -378113.5
$ws.Range("H132").Value = 1880287.5
$ws.Range("I132").Value = 3653.1765
$ws.Range("K132").Value = 10959.5295
$ws.Range("M132").Value = -8429.529500000001

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 251975
$ws.Range("I107").Value = 501200
$ws.Range("J107").Value = 2750
$ws.Range("K107").Value = 501200
$ws.Range("L107").Value = 2750
$ws.Range("M107").Value = -499280
$ws.Range("N107").Value = -6590
$ws.Range("H134").Value = 2624.5588
$ws.Range("I134").Value = 2548.84
$ws.Range("J134").Value = 2834.889
$ws.Range("K134").Value = 7646.52
$ws.Range("L134").Value = 8504.667000000001
$ws.Range("M134").Value = -5111.52
$ws.Range("N134").Value = -13574.667

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 10875578
$ws.Range("I134").Value = 15632784
$ws.Range("K134").Value = 46898352
$ws.Range("M134").Value = -46895817

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 904.3158
$ws.Range("I5").Value = 447.0909
$ws.Range("J5").Value = 1533
$ws.Range("K5").Value = 1341.2727
$ws.Range("L5").Value = 4599
$ws.Range("M5").Value = -1229.2727
$ws.Range("N5").Value = -4823
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H20").Value = 974
$ws.Range("J20").Value = 993.1818
$ws.Range("L20").Value = 2979.5454
$ws.Range("N20").Value = -3433.5454
$ws.Range("H75").Value = 7800
$ws.Range("I75").Value = 600
$ws.Range("J75").Value = 15000
$ws.Range("K75").Value = 1800
$ws.Range("L75").Value = 45000
$ws.Range("M75").Value = -802
$ws.Range("N75").Value = -46996
$ws.Range("H78").Value = 7800
$ws.Range("I78").Value = 600
$ws.Range("J78").Value = 15000
$ws.Range("K78").Value = 5400
$ws.Range("L78").Value = 135000
$ws.Range("M78").Value = -408
$ws.Range("N78").Value = -144984
$ws.Range("H92").Value = 503.14285
$ws.Range("I92").Value = 503.14285
$ws.Range("K92").Value = 1509.42855
$ws.Range("M92").Value = -261.4285500000001
$ws.Range("H126").Value = 2845.3225
$ws.Range("I126").Value = 1015
$ws.Range("J126").Value = 2971.5518
$ws.Range("K126").Value = 3045
$ws.Range("L126").Value = 8914.6554
$ws.Range("M126").Value = 1895
$ws.Range("N126").Value = -18794.6554
$ws.Range("H135").Value = 904.3158
$ws.Range("I135").Value = 447.0909
$ws.Range("J135").Value = 1533
$ws.Range("K135").Value = 4023.8181
$ws.Range("L135").Value = 13797
$ws.Range("M135").Value = -1488.8181
$ws.Range("N135").Value = -18867

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 464.4
$ws.Range("I107").Value = 474
$ws.Range("J107").Value = 450
$ws.Range("K107").Value = 474
$ws.Range("L107").Value = 450
$ws.Range("M107").Value = 1446
$ws.Range("N107").Value = -4290
$ws.Range("H132").Value = 3555.4285
$ws.Range("I132").Value = 3152.05
$ws.Range("K132").Value = 9456.150000000001
$ws.Range("M132").Value = -6926.150000000001

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 8841020
$ws.Range("I132").Value = 3728.2
$ws.Range("K132").Value = 11184.6
$ws.Range("M132").Value = -8654.599999999999